$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header cells: "Fighter A" -> "Fighter 1", "Fighter B" -> "Fighter 2"
$ws.Range("A1").Value = "Fighter 1"
$ws.Range("B1").Value = "Fighter 2"

# Match the resulting column widths (Excel auto-adjusted these after the edit)
$ws.Columns.Item(1).ColumnWidth = 17.4986979166667
$ws.Columns.Item(2).ColumnWidth = 19.7213541666667
$ws.Columns.Item(3).ColumnWidth = 18.7213541666667
$ws.Columns.Item(4).ColumnWidth = 21.1666666666667

# Match the resulting active cell selection
$null = $ws.Range("B13").Select()
